# New crime data collected - weekly NYPD 105th Precinct CompStat report update
# Updates the report header (volume number and week-covering date range)
# and refreshes crime-statistics figures for rows 15-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/12/2024  Through  8/18/2024"

# --- Crime statistics table updates (rows 15-33) ---
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 21
$ws.Range("K15").Value = 75
$ws.Range("L15").Value = 162.5
$ws.Range("M15").Value = 61.538461538461
$ws.Range("N15").Value = -19.230769230769
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 84
$ws.Range("J16").Value = 108
$ws.Range("K16").Value = -22.222222222222
$ws.Range("L16").Value = -28.813559322033
$ws.Range("M16").Value = -60.189573459715
$ws.Range("N16").Value = -87.037037037037
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -63.636363636363
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -23.529411764705
$ws.Range("I17").Value = 278
$ws.Range("J17").Value = 287
$ws.Range("K17").Value = -3.135888501742
$ws.Range("L17").Value = -2.797202797202
$ws.Range("M17").Value = 46.315789473684
$ws.Range("N17").Value = 11.646586345381
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -40.909090909090
$ws.Range("I18").Value = 118
$ws.Range("J18").Value = 140
$ws.Range("K18").Value = -15.714285714285
$ws.Range("L18").Value = 7.272727272727
$ws.Range("M18").Value = -47.787610619469
$ws.Range("N18").Value = -88.164493480441
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 84
$ws.Range("H19").Value = -40.476190476190
$ws.Range("I19").Value = 372
$ws.Range("J19").Value = 438
$ws.Range("K19").Value = -15.068493150684
$ws.Range("L19").Value = -9.927360774818
$ws.Range("M19").Value = 35.272727272727
$ws.Range("N19").Value = -1.063829787234
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 40
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = 29.032258064516
$ws.Range("I20").Value = 250
$ws.Range("J20").Value = 205
$ws.Range("K20").Value = 21.951219512195
$ws.Range("L20").Value = 65.562913907284
$ws.Range("M20").Value = 4.166666666666
$ws.Range("N20").Value = -88.312295465170
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = -27.659574468085
$ws.Range("F21").Value = 144
$ws.Range("G21").Value = 190
$ws.Range("H21").Value = -24.210526315789
$ws.Range("I21").Value = 1125
$ws.Range("J21").Value = 1192
$ws.Range("K21").Value = -5.620805369127
$ws.Range("L21").Value = 3.305785123966
$ws.Range("M21").Value = -3.350515463917
$ws.Range("N21").Value = -74.702046323364
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -27.777777777777
$ws.Range("F24").Value = 87
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -23.684210526315
$ws.Range("I24").Value = 719
$ws.Range("J24").Value = 784
$ws.Range("K24").Value = -8.290816326530
$ws.Range("L24").Value = -23.673036093418
$ws.Range("M24").Value = 35.660377358490
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 22.222222222222
$ws.Range("I25").Value = 169
$ws.Range("J25").Value = 132
$ws.Range("K25").Value = 28.030303030303
$ws.Range("L25").Value = -26.521739130434
$ws.Range("C26").Value = 24
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 166.666666666667
$ws.Range("F26").Value = 66
$ws.Range("G26").Value = 59
$ws.Range("H26").Value = 11.864406779661
$ws.Range("I26").Value = 502
$ws.Range("J26").Value = 404
$ws.Range("K26").Value = 24.257425742574
$ws.Range("L26").Value = 19.523809523809
$ws.Range("M26").Value = 23.341523341523
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 400
$ws.Range("I27").Value = 32
$ws.Range("K27").Value = 45.454545454545
$ws.Range("L27").Value = 60
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 35
$ws.Range("K28").Value = 16.666666666666
$ws.Range("L28").Value = -5.405405405405
$ws.Range("F29").Value = 1
$ws.Range("L29").Value = -38.888888888888
$ws.Range("F30").Value = 1
$ws.Range("L30").Value = -28.571428571428
$ws.Range("G33").Value = 3
$ws.Range("J33").Value = 7
$ws.Range("K33").Value = -85.714285714285